$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.073.53"
$ws.Range("D3").Value = "1.681.13"
$ws.Range("E3").Value = "  +0.77%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.81"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("E6").Value = "  -3.26%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0623"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.28"
$ws.Range("E10").Value = "  +5.01%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.918.60"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "1.687.27"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.11"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.535"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.42"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "27.085.11"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.15"
$ws.Range("E18").Value = "  +2.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "235.92"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").Value = "  +1.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.26"
$ws.Range("E23").Value = "  +1.25%  "
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.49"
$ws.Range("E25").Value = "  +0.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.28"
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.53"
$ws.Range("E27").Value = "  +3.80%  "
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.37"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "1.548.95"
$ws.Range("E33").Value = "  +6.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  +4.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.585"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.914"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +2.37%  "
$ws.Range("E40").Value = "  +7.51%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.88"
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.25"
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").Value = "1.823.68"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.50"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  +1.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.01"
$ws.Range("E51").Value = "  +5.83%  "
